# Updates cryptocurrency price/volume data in the worksheet to reflect the
# latest scrape (GitHub Actions cryptos list refresh).
# Target cells are written as literal text (not auto-parsed numbers) to
# preserve the source formatting exactly, matching how the sheet stores
# prices/percentages as inline strings. A leading apostrophe is used for
# values that would otherwise be auto-converted to numbers by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.124.16"
$ws.Range("E2").Value = "  -2.02%  "
$ws.Range("D3").Value = "2.742.28"
$ws.Range("E3").Value = "  -2.62%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'352.15"
$ws.Range("E5").Value = "  -2.64%  "
$ws.Range("D6").Value = "'106.64"
$ws.Range("E6").Value = "  -3.86%  "
$ws.Range("D7").Value = "'0.546"
$ws.Range("E7").Value = "  -3.39%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").Value = "'0.575"
$ws.Range("E9").Value = "  -4.09%  "
$ws.Range("D10").Value = "'38.82"
$ws.Range("E10").Value = "  -4.11%  "
$ws.Range("E11").Value = "  +2.81%  "
$ws.Range("D12").Value = "'0.0829"
$ws.Range("E12").Value = "  -3.62%  "
$ws.Range("D13").Value = "'19.56"
$ws.Range("E13").Value = "  -1.18%  "
$ws.Range("D14").Value = "'7.44"
$ws.Range("E14").Value = "  -4.37%  "
$ws.Range("D15").Value = "3.188.79"
$ws.Range("E15").Value = "  -2.18%  "
$ws.Range("D16").Value = "2.754.63"
$ws.Range("E16").Value = "  -3.18%  "
$ws.Range("D17").Value = "'0.914"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("D18").Value = "51.141.60"
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("D19").Value = "'7.55"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("D20").Value = "'3.01"
$ws.Range("E20").Value = "  -3.51%  "
$ws.Range("D21").Value = "'12.86"
$ws.Range("E21").Value = "  -3.32%  "
$ws.Range("D22").Value = "0.0₃0953"
$ws.Range("E22").Value = "  -4.04%  "
$ws.Range("D23").Value = "'69.14"
$ws.Range("E23").Value = "  -1.41%  "
$ws.Range("D24").Value = "'262.57"
$ws.Range("E24").Value = "  -3.71%  "
$ws.Range("D25").Value = "'2.71"
$ws.Range("E25").Value = "  -3.46%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").Value = "'25.80"
$ws.Range("E27").Value = "  -3.43%  "
$ws.Range("E28").Value = "  +12.45%  "
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("D30").Value = "'10.00"
$ws.Range("E30").Value = "  -2.38%  "
$ws.Range("B31").Value = "OKB"
$ws.Range("C31").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D31").Value = "'51.51"
$ws.Range("E31").Value = "  -1.21%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'6.01"
$ws.Range("E32").Value = "  +3.06%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").Value = "'34.23"
$ws.Range("E33").Value = "  -0.65%  "
$ws.Range("D34").Value = "'0.0438"
$ws.Range("E34").Value = "  -7.60%  "
$ws.Range("D35").Value = "'0.0824"
$ws.Range("E35").Value = "  -2.44%  "
$ws.Range("D36").Value = "'5.10"
$ws.Range("E36").Value = "  -6.90%  "
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("D38").Value = "'18.20"
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("D39").Value = "'3.11"
$ws.Range("E39").Value = "  -3.51%  "
$ws.Range("D40").Value = "'1.92"
$ws.Range("E40").Value = "  -4.85%  "
$ws.Range("E41").Value = "  -3.31%  "
$ws.Range("D42").Value = "'2.45"
$ws.Range("E42").Value = "  -3.64%  "
$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("D44").Value = "'119.56"
$ws.Range("E44").Value = "  -4.31%  "
$ws.Range("D45").Value = "'21.63"
$ws.Range("E45").Value = "  -3.92%  "
$ws.Range("D46").Value = "2.066.16"
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("D47").Value = "'2.32"
$ws.Range("E47").Value = "  -0.92%  "
$ws.Range("D48").Value = "'3.19"
$ws.Range("E48").Value = "  -2.96%  "
$ws.Range("D49").Value = "'0.909"
$ws.Range("E49").Value = "  -4.44%  "
$ws.Range("D50").Value = "'5.43"
$ws.Range("E50").Value = "  -7.03%  "
$ws.Range("D51").Value = "'58.52"
$ws.Range("E51").Value = "  -2.53%  "
